$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.014.78"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.58%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.448.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.13%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.06%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'585.36"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.07%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'174.09"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.57%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.04%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.604"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +1.06%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'3.448.15"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.12%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.01%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'6.97"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.87%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.414"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.77%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'4.047.29"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.06%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.93%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'29.17"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -5.05%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'65.983.73"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.64%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.08%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'3.444.75"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.22%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'5.96"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.40%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'13.84"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.38%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'370.73"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -1.18%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'7.60"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.79%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'72.29"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.37%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.09%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.15%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +4.63%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'9.74"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.59%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +3.68%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.00%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.01%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.58%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'1.99"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.47%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.05%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -4.47%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'7.02"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.04%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +1.60%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'161.61"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +1.55%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.881"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.62%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'28.40"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +5.24%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'1.79"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.84%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.27%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.780.74"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +3.28%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.87%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.42%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0687"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.64%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'25.00"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.60%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'39.91"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.68%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.0292"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.27%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'329.05"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.90%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.51%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'6.27"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.74%  "
$ws.Range('E51').Style = 'Normal'
